$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7145.077
$ws.Range("I62").Value = 4582.3335
$ws.Range("J62").Value = 9341.714
$ws.Range("K62").Value = 4582.3335
$ws.Range("L62").Value = 9341.714
$ws.Range("M62").Value = -3958.3335
$ws.Range("N62").Value = -10589.714

$ws.Range("H65").Value = 7145.077
$ws.Range("I65").Value = 4582.3335
$ws.Range("J65").Value = 9341.714
$ws.Range("K65").Value = 22911.6675
$ws.Range("L65").Value = 46708.57
$ws.Range("M65").Value = -19791.6675
$ws.Range("N65").Value = -52948.57

$ws.Range("H116").Value = 4116
$ws.Range("I116").Value = 2164.3635
$ws.Range("K116").Value = 2164.3635
$ws.Range("M116").Value = 1277.6365

$ws.Range("H129").Value = 162189.5
$ws.Range("I129").Value = 219.77777
$ws.Range("K129").Value = 659.33331
$ws.Range("M129").Value = 4340.66669

$ws.Range("H137").Value = 1187.5714
$ws.Range("I137").Value = 1199.919
$ws.Range("K137").Value = 3599.757000000001
$ws.Range("M137").Value = -1049.757000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6185.573
$ws.Range("I32").Value = 4219.984
$ws.Range("K32").Value = 4219.984
$ws.Range("M32").Value = -3932.984

$ws.Range("H45").Value = 2679.0952
$ws.Range("I45").Value = 1858.0834
$ws.Range("J45").Value = 3773.7778
$ws.Range("K45").Value = 1858.0834
$ws.Range("L45").Value = 3773.7778
$ws.Range("M45").Value = -1481.0834
$ws.Range("N45").Value = -4527.7778

$ws.Range("H61").Value = 2298.5454
$ws.Range("I61").Value = 1660.5
$ws.Range("K61").Value = 1660.5
$ws.Range("M61").Value = -1448.5

$ws.Range("H97").Value = 34483572
$ws.Range("I97").Value = 425.45
$ws.Range("K97").Value = 425.45
$ws.Range("M97").Value = 70.55000000000001

$ws.Range("H122").Value = 1646.9375
$ws.Range("J122").Value = 5344.5
$ws.Range("L122").Value = 16033.5
$ws.Range("N122").Value = -20933.5

$ws.Range("H132").Value = 12631.718
$ws.Range("I132").Value = 1630.8422
$ws.Range("K132").Value = 4892.5266
$ws.Range("M132").Value = -2362.5266

$ws.Range("H136").Value = 2298.5454
$ws.Range("I136").Value = 1660.5
$ws.Range("K136").Value = 4981.5
$ws.Range("M136").Value = -2431.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 751.4286
$ws.Range("I94").Value = 625.3333
$ws.Range("J94").Value = 1066.6666
$ws.Range("K94").Value = 625.3333
$ws.Range("L94").Value = 1066.6666
$ws.Range("M94").Value = -174.3333
$ws.Range("N94").Value = -1968.6666

$ws.Range("H105").Value = 1769.6
$ws.Range("I105").Value = 1506.6666
$ws.Range("J105").Value = 1852.6316
$ws.Range("K105").Value = 1506.6666
$ws.Range("L105").Value = 1852.6316
$ws.Range("M105").Value = 240.3334
$ws.Range("N105").Value = -5346.6316

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 451.375
$ws.Range("J22").Value = 550.1667
$ws.Range("L22").Value = 550.1667
$ws.Range("N22").Value = -1250.1667

$ws.Range("H31").Value = 4160.552
$ws.Range("I31").Value = 4138.5557
$ws.Range("J31").Value = 4170.45
$ws.Range("K31").Value = 4138.5557
$ws.Range("L31").Value = 4170.45
$ws.Range("M31").Value = -3843.5557
$ws.Range("N31").Value = -4760.45

$ws.Range("H34").Value = 4160.552
$ws.Range("I34").Value = 4138.5557
$ws.Range("J34").Value = 4170.45
$ws.Range("K34").Value = 4138.5557
$ws.Range("L34").Value = 4170.45
$ws.Range("M34").Value = -3936.5557
$ws.Range("N34").Value = -4574.45

$ws.Range("H58").Value = 20798.076
$ws.Range("I58").Value = 1577.7858
$ws.Range("J58").Value = 43221.75
$ws.Range("K58").Value = 1577.7858
$ws.Range("L58").Value = 43221.75
$ws.Range("M58").Value = -1374.7858
$ws.Range("N58").Value = -43627.75

$ws.Range("H134").Value = 1392.8334
$ws.Range("I134").Value = 1300
$ws.Range("J134").Value = 1522.8
$ws.Range("K134").Value = 3900
$ws.Range("L134").Value = 4568.4
$ws.Range("M134").Value = -1365
$ws.Range("N134").Value = -9638.4

$ws.Range("H136").Value = 20798.076
$ws.Range("I136").Value = 1577.7858
$ws.Range("J136").Value = 43221.75
$ws.Range("K136").Value = 4733.357400000001
$ws.Range("L136").Value = 129665.25
$ws.Range("M136").Value = -2183.357400000001
$ws.Range("N136").Value = -134765.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1174.091
$ws.Range("I5").Value = 538.8
$ws.Range("J5").Value = 1703.5
$ws.Range("K5").Value = 1616.4
$ws.Range("L5").Value = 5110.5
$ws.Range("M5").Value = -1504.4
$ws.Range("N5").Value = -5334.5

$ws.Range("H55").Value = 2488.6365
$ws.Range("I55").Value = 1400
$ws.Range("J55").Value = 2597.5
$ws.Range("K55").Value = 4200
$ws.Range("L55").Value = 7792.5
$ws.Range("M55").Value = -4023
$ws.Range("N55").Value = -8146.5

$ws.Range("H80").Value = 14699.875
$ws.Range("J80").Value = 18916.666
$ws.Range("L80").Value = 56749.99800000001
$ws.Range("N80").Value = -58621.99800000001

$ws.Range("H83").Value = 14699.875
$ws.Range("J83").Value = 18916.666
$ws.Range("L83").Value = 170249.994
$ws.Range("N83").Value = -179609.994

$ws.Range("H86").Value = 83334140
$ws.Range("I86").Value = 950.6667
$ws.Range("J86").Value = 166667330
$ws.Range("K86").Value = 2852.0001
$ws.Range("L86").Value = 500001990
$ws.Range("M86").Value = -1666.0001
$ws.Range("N86").Value = -500004362

$ws.Range("H89").Value = 83334140
$ws.Range("I89").Value = 950.6667
$ws.Range("J89").Value = 166667330
$ws.Range("K89").Value = 8556.0003
$ws.Range("L89").Value = 1500005970
$ws.Range("M89").Value = -2628.0003
$ws.Range("N89").Value = -1500017826

$ws.Range("H118").Value = 2571.5
$ws.Range("I118").Value = 1357.25
$ws.Range("J118").Value = 5000
$ws.Range("K118").Value = 4071.75
$ws.Range("L118").Value = 15000
$ws.Range("M118").Value = -2828.75
$ws.Range("N118").Value = -17486

$ws.Range("H122").Value = 757.2143
$ws.Range("I122").Value = 435
$ws.Range("J122").Value = 810.9167
$ws.Range("K122").Value = 3915
$ws.Range("L122").Value = 7298.2503
$ws.Range("M122").Value = -1465
$ws.Range("N122").Value = -12198.2503

$ws.Range("H130").Value = 1400
$ws.Range("I130").Value = 1400
$ws.Range("K130").Value = 4200
$ws.Range("M130").Value = 820

$ws.Range("H131").Value = 692.46
$ws.Range("I131").Value = 300
$ws.Range("J131").Value = 696.42426
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 2089.27278
$ws.Range("M131").Value = 4140
$ws.Range("N131").Value = -12169.27278

$ws.Range("H132").Value = 541.4286
$ws.Range("J132").Value = 500
$ws.Range("L132").Value = 4500
$ws.Range("N132").Value = -9560

$ws.Range("H133").Value = 4166
$ws.Range("I133").Value = 1129.6666
$ws.Range("J133").Value = 4994.091
$ws.Range("K133").Value = 3388.9998
$ws.Range("L133").Value = 14982.273
$ws.Range("M133").Value = 1671.0002
$ws.Range("N133").Value = -25102.273

$ws.Range("H134").Value = 2265
$ws.Range("I134").Value = 1556.5416
$ws.Range("J134").Value = 4694
$ws.Range("K134").Value = 4669.6248
$ws.Range("L134").Value = 14082
$ws.Range("M134").Value = 400.3752000000004
$ws.Range("N134").Value = -24222

$ws.Range("H135").Value = 1174.091
$ws.Range("I135").Value = 538.8
$ws.Range("J135").Value = 1703.5
$ws.Range("K135").Value = 4849.2
$ws.Range("L135").Value = 15331.5
$ws.Range("M135").Value = -2314.2
$ws.Range("N135").Value = -20401.5

$ws.Range("H138").Value = 2081.0476
$ws.Range("I138").Value = 1821.7059
$ws.Range("J138").Value = 3183.25
$ws.Range("K138").Value = 5465.1177
$ws.Range("L138").Value = 9549.75
$ws.Range("M138").Value = -325.1176999999998
$ws.Range("N138").Value = -19829.75

$ws.Range("H139").Value = 1963.8368
$ws.Range("I139").Value = 1034.7693
$ws.Range("J139").Value = 3014.087
$ws.Range("K139").Value = 3104.3079
$ws.Range("L139").Value = 9042.261
$ws.Range("M139").Value = 2035.6921
$ws.Range("N139").Value = -19322.261

$ws.Range("H140").Value = 1957.56
$ws.Range("I140").Value = 1246.25
$ws.Range("J140").Value = 3222.111
$ws.Range("K140").Value = 3738.75
$ws.Range("L140").Value = 9666.332999999999
$ws.Range("M140").Value = 1441.25
$ws.Range("N140").Value = -20026.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1626.875
$ws.Range("I97").Value = 1433.0769
$ws.Range("J97").Value = 2466.6667
$ws.Range("K97").Value = 1433.0769
$ws.Range("L97").Value = 2466.6667
$ws.Range("M97").Value = -937.0769
$ws.Range("N97").Value = -3458.6667

$ws.Range("H113").Value = 3847.8823
$ws.Range("I113").Value = 4857.4346
$ws.Range("J113").Value = 1737
$ws.Range("K113").Value = 4857.4346
$ws.Range("L113").Value = 1737
$ws.Range("M113").Value = -2687.4346
$ws.Range("N113").Value = -6077

$ws.Range("H139").Value = 25174.615
$ws.Range("J139").Value = 25174.615
$ws.Range("L139").Value = 25174.615
$ws.Range("N139").Value = -35454.61500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 932.4737
$ws.Range("I46").Value = 931.75
$ws.Range("K46").Value = 931.75
$ws.Range("M46").Value = -743.75

$ws.Range("H140").Value = 43717.4
$ws.Range("J140").Value = 43717.4
$ws.Range("L140").Value = 43717.4
$ws.Range("N140").Value = -54077.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31872

$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99360

$ws.Range("H138").Value = 49453
$ws.Range("J138").Value = 49453
$ws.Range("L138").Value = 49453
$ws.Range("N138").Value = -59733
